$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value2 = "29.494.69"
$ws.Range("E2").Value2 = "  +0.94%  "

# Row 3
$ws.Range("D3").Value2 = "1.882.98"
$ws.Range("E3").Value2 = "  +1.61%  "

# Row 4
$ws.Range("D4").Value2 = "'0.9992"
$ws.Range("E4").Value2 = "  -0.22%  "

# Row 5
$ws.Range("D5").Value2 = "'0.7165"
$ws.Range("E5").Value2 = "  +2.31%  "

# Row 6
$ws.Range("D6").Value2 = "'242.58"
$ws.Range("E6").Value2 = "  +2.10%  "

# Row 7
$ws.Range("D7").Value2 = "'0.9995"
$ws.Range("E7").Value2 = "  -0.21%  "

# Row 8
$ws.Range("D8").Value2 = "'0.07927"
$ws.Range("E8").Value2 = "  +0.37%  "

# Row 9
$ws.Range("E9").Value2 = "  +3.64%  "

# Row 10
$ws.Range("D10").Value2 = "'25.40"
$ws.Range("E10").Value2 = "  +7.89%  "

# Row 11
$ws.Range("D11").Value2 = "'0.08283"
$ws.Range("E11").Value2 = "  +1.29%  "

# Row 12
$ws.Range("E12").Value2 = "  +4.11%  "

# Row 13
$ws.Range("D13").Value2 = "1.878.33"
$ws.Range("E13").Value2 = "  +1.45%  "

# Row 14
$ws.Range("D14").Value2 = "'5.299"
$ws.Range("E14").Value2 = "  +2.34%  "

# Row 15
$ws.Range("D15").Value2 = "'91.55"
$ws.Range("E15").Value2 = "  +2.36%  "

# Row 16
$ws.Range("D16").Value2 = "29.496.34"
$ws.Range("E16").Value2 = "  +0.95%  "

# Row 17
$ws.Range("D17").Value2 = "'5.960"
$ws.Range("E17").Value2 = "  +2.69%  "

# Row 18
$ws.Range("D18").Value2 = "'248.19"
$ws.Range("E18").Value2 = "  +5.14%  "

# Row 19
$ws.Range("D19").Value2 = "'0.000007887"
$ws.Range("E19").Value2 = "  +1.03%  "

# Row 20
$ws.Range("D20").Value2 = "'13.39"
$ws.Range("E20").Value2 = "  +1.67%  "

# Row 21
$ws.Range("B21").Value2 = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value2 = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value2 = "2.142.09"
$ws.Range("E21").Value2 = "  +2.38%  "

# Row 22
$ws.Range("B22").Value2 = "Dai"
$ws.Range("C22").Value2 = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value2 = "'0.9992"
$ws.Range("E22").Value2 = "  -0.14%  "

# Row 23
$ws.Range("B23").Value2 = "Chainlink"
$ws.Range("C23").Value2 = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").Value2 = "'7.991"
$ws.Range("E23").Value2 = "  +6.66%  "

# Row 24
$ws.Range("B24").Value2 = "BinanceUSD"
$ws.Range("C24").Value2 = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D24").Value2 = "'0.9996"
$ws.Range("E24").Value2 = "  -0.20%  "

# Row 25
$ws.Range("B25").Value2 = "Stellar"
$ws.Range("C25").Value2 = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D25").Value2 = "'0.1614"
$ws.Range("E25").Value2 = "  +14.24%  "

# Row 26
$ws.Range("B26").Value2 = "Monero"
$ws.Range("C26").Value2 = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value2 = "'163.43"
$ws.Range("E26").Value2 = "  +0.42%  "

# Row 27
$ws.Range("B27").Value2 = "Cosmos"
$ws.Range("C27").Value2 = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value2 = "'9.076"
$ws.Range("E27").Value2 = "  +2.63%  "

# Row 28
$ws.Range("B28").Value2 = "EthereumClassic"
$ws.Range("C28").Value2 = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value2 = "'18.39"
$ws.Range("E28").Value2 = "  +1.97%  "

# Row 29
$ws.Range("B29").Value2 = "Toncoin"
$ws.Range("C29").Value2 = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value2 = "'1.357"
$ws.Range("E29").Value2 = "  -3.34%  "

# Row 30
$ws.Range("B30").Value2 = "PancakeSwap"
$ws.Range("C30").Value2 = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value2 = "'1.504"
$ws.Range("E30").Value2 = "  +2.36%  "

# Row 31
$ws.Range("B31").Value2 = "Filecoin"
$ws.Range("C31").Value2 = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value2 = "'4.405"
$ws.Range("E31").Value2 = "  +2.12%  "

# Row 32
$ws.Range("B32").Value2 = "InternetComputer(DFINITY)"
$ws.Range("C32").Value2 = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value2 = "'4.127"
$ws.Range("E32").Value2 = "  +2.89%  "

# Row 33
$ws.Range("B33").Value2 = "Hedera"
$ws.Range("C33").Value2 = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value2 = "'0.05298"
$ws.Range("E33").Value2 = "  +3.00%  "

# Row 34
$ws.Range("B34").Value2 = "LidoDAOToken"
$ws.Range("C34").Value2 = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").Value2 = "'1.952"
$ws.Range("E34").Value2 = "  +2.34%  "

# Row 35
$ws.Range("B35").Value2 = "ARBITRUM"
$ws.Range("C35").Value2 = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value2 = "'1.203"
$ws.Range("E35").Value2 = "  +3.60%  "

# Row 36
$ws.Range("B36").Value2 = "ImmutableX"
$ws.Range("C36").Value2 = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value2 = "'0.7297"
$ws.Range("E36").Value2 = "  +3.11%  "

# Row 37
$ws.Range("B37").Value2 = "HuobiToken"
$ws.Range("C37").Value2 = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value2 = "'2.675"
$ws.Range("E37").Value2 = "  -0.30%  "

# Row 38
$ws.Range("B38").Value2 = "VeChain"
$ws.Range("C38").Value2 = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value2 = "'0.01875"
$ws.Range("E38").Value2 = "  +1.69%  "

# Row 39
$ws.Range("B39").Value2 = "Maker"
$ws.Range("C39").Value2 = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value2 = "1.228.31"
$ws.Range("E39").Value2 = "  +6.39%  "

# Row 40
$ws.Range("B40").Value2 = "MXToken"
$ws.Range("C40").Value2 = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value2 = "'2.736"
$ws.Range("E40").Value2 = "  +1.05%  "

# Row 41
$ws.Range("B41").Value2 = "TrustWalletToken"
$ws.Range("C41").Value2 = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value2 = "'0.9148"
$ws.Range("E41").Value2 = "  -1.37%  "

# Row 42
$ws.Range("B42").Value2 = "Aave"
$ws.Range("C42").Value2 = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value2 = "'74.95"
$ws.Range("E42").Value2 = "  +7.02%  "

# Row 43
$ws.Range("B43").Value2 = "FraxShare"
$ws.Range("C43").Value2 = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value2 = "'6.208"
$ws.Range("E43").Value2 = "  +3.96%  "

# Row 44
$ws.Range("B44").Value2 = "PaxDollar"
$ws.Range("C44").Value2 = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").Value2 = "'0.9996"
$ws.Range("E44").Value2 = "  -0.18%  "

# Row 45
$ws.Range("B45").Value2 = "Quant"
$ws.Range("C45").Value2 = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value2 = "'102.74"
$ws.Range("E45").Value2 = "  +0.22%  "

# Row 46
$ws.Range("B46").Value2 = "RocketPoolETH"
$ws.Range("C46").Value2 = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value2 = "2.040.26"
$ws.Range("E46").Value2 = "  +2.54%  "

# Row 47
$ws.Range("B47").Value2 = "Mantle"
$ws.Range("C47").Value2 = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value2 = "'0.5277"
$ws.Range("E47").Value2 = "  -0.29%  "

# Row 48
$ws.Range("B48").Value2 = "RenderToken"
$ws.Range("C48").Value2 = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value2 = "'1.789"
$ws.Range("E48").Value2 = "  +3.21%  "

# Row 49
$ws.Range("B49").Value2 = "SynthetixNetwork"
$ws.Range("C49").Value2 = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D49").Value2 = "'2.936"
$ws.Range("E49").Value2 = "  +12.19%  "

# Row 50
$ws.Range("B50").Value2 = "EnergySwap"
$ws.Range("C50").Value2 = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value2 = "'9.344"
$ws.Range("E50").Value2 = "  +2.63%  "

# Row 51
$ws.Range("B51").Value2 = "TheSandbox"
$ws.Range("C51").Value2 = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D51").Value2 = "'0.4345"
$ws.Range("E51").Value2 = "  +2.52%  "
